$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Medstar POB (North -> South Tower), updated year built & area ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319

# --- Row 3: 1801 Pennsylvania Ave. -> 1801 Pennsylvania Avenue, LLC ---
$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

# --- Row 4: GSA 300 E Street SW - address fix & new owner & area ---
$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

# --- Row 5: Paul H.Nitze - only area (GHG/EUI) updated ---
$ws.Range("L5").Value = 58717

# --- Row 6: President Madison Apartments -> Hampton House (replaced property) ---
$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580

# --- Row 7: 3303 Water Street - postal code & area updated ---
$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

# --- Row 8: 15th and H Street Associates LLP - address fix ---
$ws.Range("E8").Value = "1428 H ST NW"

# --- Row 9: Eastern Market - no data changes ---

# --- Row 10: DPW Vehicle Maintenance Facility 2 -> School Without Walls @ Francis Stevens (replaced property) ---
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991

# Year Built column no longer carries the date-style formatting (s="2") in any row
$ws.Range("I2:I10").Style = "Normal"

# Selection moved to the full data block instead of the old Owner-email column
$ws.Range("A1:L10").Select()
